$d = $word.ActiveDocument

# --- 1) "Programa" section, Portuguese paragraph: split the single run
#        into four sentences joined by manual line breaks (<w:br/>). ---
$old1 = "- Revisão das regras de segurança laboratorial - Introdução à análise qualitativa: Definições, objetivos e limitações. - Identificação dos cátions do grupo I (K+, Na+ e NH4+); grupo II (Mg2+, Ca2+ e Ba2+); grupo III (Al3+, Fe3+, Mn2+).- Estudo dos ânions e suas aplicações em análises ambientais (Cl e suas espécies, SO42-, CO32-, S2-, NO3-)."
$new1 = "- Revisão das regras de segurança laboratorial ^l- Introdução à análise qualitativa: Definições, objetivos e limitações. ^l- Identificação dos cátions do grupo I (K+, Na+ e NH4+); grupo II (Mg2+, Ca2+ e Ba2+); grupo III (Al3+, Fe3+, Mn2+).^l- Estudo dos ânions e suas aplicações em análises ambientais (Cl e suas espécies, SO42-, CO32-, S2-, NO3-)."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# --- 2) "Programa" section, English (italic) paragraph: same kind of split. ---
$old2 = "- Review of laboratory safety rules- Introduction to qualitative analysis: Definitions, objectives and limitations.- Identification of group I cations (K+, Na+ and NH4+); group II (Mg2+, Ca2+ and Ba2+); group III (Al3+, Fe3+, Fe2+, Mn2+).- Study of anions and their applications in environmental analysis (Cl and its species, SO42-, CO32-, S2-, NO3-)."
$new2 = "- Review of laboratory safety rules^l- Introduction to qualitative analysis: Definitions, objectives and limitations.^l- Identification of group I cations (K+, Na+ and NH4+); group II (Mg2+, Ca2+ and Ba2+); group III (Al3+, Fe3+, Fe2+, Mn2+).^l- Study of anions and their applications in environmental analysis (Cl and its species, SO42-, CO32-, S2-, NO3-)."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# --- 3) "Bibliografia" section: split the run into one line per reference. ---
$old3 = "Baird, C., Michael C.C.  Environmental chemistry. Editora: New York: Freeman, 5a edição, 2012Baird, C.,Michael C.C.  Química ambiental. Editora: Porto Alegre: Bookman,  4a edição, 2011Harris, D. C. Análise Química Quantitativa. Editora: LTC, 8a edição, 2012Krug, F.J., Rocha F.R.P. Métodos de preparo de amostras para análise elementar. Editora EditSBQ, 1a edição, 2016Luna, A. Química analítica ambiental. Editora: EdUERJ, 1a edição, 2003Rocha, J.C., Rosa, A.H., Cardoso, A.A. Introdução à química ambiental. Editora: Porto Alegre: Bookman, 2a edição, 2009.Skoog, D. A, West, D. M., Holler, F. J., Crouch, S. R. Fundamentos de Química Analítica. Editora: Thomson, tradução da 8ª edição, 2006"
$new3 = "Baird, C., Michael C.C.  Environmental chemistry. Editora: New York: Freeman, 5a edição, 2012^lBaird, C.,Michael C.C.  Química ambiental. Editora: Porto Alegre: Bookman,  4a edição, 2011^lHarris, D. C. Análise Química Quantitativa. Editora: LTC, 8a edição, 2012^lKrug, F.J., Rocha F.R.P. Métodos de preparo de amostras para análise elementar. Editora EditSBQ, 1a edição, 2016^lLuna, A. Química analítica ambiental. Editora: EdUERJ, 1a edição, 2003^lRocha, J.C., Rosa, A.H., Cardoso, A.A. Introdução à química ambiental. Editora: Porto Alegre: Bookman, 2a edição, 2009.^lSkoog, D. A, West, D. M., Holler, F. J., Crouch, S. R. Fundamentos de Química Analítica. Editora: Thomson, tradução da 8ª edição, 2006"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
